# Auto-generated script applying cryptos list update diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D contains numeric-looking text values that must remain text,
# not be auto-converted to numbers by Excel. Force text format first,
# write the values, then clear the temporary format so cell styling
# matches the original (unstyled) cells.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "95.928.18"
$ws.Range("E2").Value = "  -0.73%  "
$ws.Range("D3").Value = "3.669.83"
$ws.Range("E3").Value = "  +1.05%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "235.29"
$ws.Range("E5").Value = "  -2.80%  "
$ws.Range("D6").Value = "1.87"
$ws.Range("E6").Value = "  +4.35%  "
$ws.Range("D7").Value = "648.57"
$ws.Range("E7").Value = "  -0.79%  "
$ws.Range("D8").Value = "0.420"
$ws.Range("E8").Value = "  +0.49%  "
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("D10").Value = "1.04"
$ws.Range("E10").Value = "  -2.07%  "
$ws.Range("D11").Value = "3.668.17"
$ws.Range("E11").Value = "  +1.07%  "
$ws.Range("D12").Value = "44.12"
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("D13").Value = "0.203"
$ws.Range("E13").Value = "  -0.22%  "
$ws.Range("D14").Value = "6.69"
$ws.Range("E14").Value = "  +2.37%  "
$ws.Range("D15").Value = "0.0000284"
$ws.Range("E15").Value = "  +10.40%  "
$ws.Range("D16").Value = "4.353.53"
$ws.Range("E16").Value = "  +1.08%  "
$ws.Range("D17").Value = "96.135.72"
$ws.Range("E17").Value = "  -0.26%  "
$ws.Range("D18").Value = "3.680.61"
$ws.Range("E18").Value = "  +2.10%  "
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").Value = "8.03"
$ws.Range("E19").Value = "  -7.32%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "12.75"
$ws.Range("E20").Value = "  -2.55%  "
$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").Value = "18.56"
$ws.Range("E21").Value = "  +0.72%  "
$ws.Range("D22").Value = "0.498"
$ws.Range("E22").Value = "  -6.14%  "
$ws.Range("D23").Value = "512.44"
$ws.Range("E23").Value = "  +0.35%  "
$ws.Range("D24").Value = "3.35"
$ws.Range("E24").Value = "  -2.44%  "
$ws.Range("D25").Value = "0.0000202"
$ws.Range("E25").Value = "  -1.04%  "
$ws.Range("D26").Value = "6.89"
$ws.Range("E26").Value = "  +0.87%  "
$ws.Range("D27").Value = "99.90"
$ws.Range("E27").Value = "  -1.01%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "3.869.66"
$ws.Range("E28").Value = "  +1.19%  "
$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").Value = "13.00"
$ws.Range("E29").Value = "  -0.87%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").Value = "0.166"
$ws.Range("E30").Value = "  -1.35%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "2.98"
$ws.Range("E31").Value = "  -2.17%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "12.00"
$ws.Range("E32").Value = "  +0.77%  "
$ws.Range("B33").Value = "Dai"
$ws.Range("C33").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D33").Value = "0.998"
$ws.Range("E33").Value = "  -0.24%  "
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").Value = "1.84"
$ws.Range("E34").Value = "  +6.56%  "
$ws.Range("B35").Value = "Cronos"
$ws.Range("C35").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D35").Value = "0.182"
$ws.Range("E35").Value = "  -1.53%  "
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D37").Value = "649.68"
$ws.Range("E37").Value = "  +5.88%  "
$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").Value = "31.89"
$ws.Range("E38").Value = "  -3.43%  "
$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D39").Value = "0.582"
$ws.Range("E39").Value = "  +1.04%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D40").Value = "8.69"
$ws.Range("E40").Value = "  -1.37%  "
$ws.Range("B41").Value = "USDe"
$ws.Range("C41").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").Value = "40.71"
$ws.Range("E42").Value = "  +5.01%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "6.71"
$ws.Range("E43").Value = "  +8.05%  "
$ws.Range("D44").Value = "2.00"
$ws.Range("E44").Value = "  +3.98%  "
$ws.Range("B45").Value = "Kaspa"
$ws.Range("C45").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D45").Value = "0.158"
$ws.Range("E45").Value = "  +1.80%  "
$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D46").Value = "0.950"
$ws.Range("E46").Value = "  -0.24%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "0.0442"
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Value = "0.428"
$ws.Range("E48").Value = "  +5.09%  "
$ws.Range("D49").Value = "23.53"
$ws.Range("E49").Value = "  -0.24%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "2.24"
$ws.Range("E50").Value = "  -2.24%  "
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").Value = "8.35"
$ws.Range("E51").Value = "  -2.41%  "

# Restore original (default) formatting on column D now that values are set
$dRange.ClearFormats()

